# Weekly data update: insert a new observation row for "Agrícola del Norte
# S.A. de Arica - Acelga" right after the most recent existing row (60),
# pushing the previously-existing rows 61-88 down to 62-89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 61; rows 61:88 shift down to 62:89, and the used
# range grows from A1:R88 to A1:R89.
$ws.Rows("61:61").Insert()

# Populate the newly inserted row 61 with this week's observation.
$ws.Range("A61").Value = 1
$ws.Range("B61").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C61").Value = "Arica y Parinacota"
$ws.Range("D61").Value = 45029
$ws.Range("E61").Value = 15
$ws.Range("F61").Value = 100112009
$ws.Range("G61").Value = "Acelga"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Segunda"
$ws.Range("J61").Value = 350
$ws.Range("K61").Value = 1800
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = 1886
$ws.Range("N61").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O61").Value = "Región de Arica y Parinacota"
$ws.Range("P61").Value = 629
$ws.Range("Q61").Value = 3
$ws.Range("R61").Value = "Hortaliza"
